$wb = $excel.ActiveWorkbook
$dataWs = $wb.Worksheets.Item("data")

# --- Update query timestamps on the "data" sheet (F2:F4) --------------------
$dataWs.Range("F2").Value = "2021-10-05 14:34:05.684081"
$dataWs.Range("F3").Value = "2021-10-05 14:34:05.684090"
$dataWs.Range("F4").Value = "2021-10-05 14:34:05.684093"

# --- Add the new "metadata" worksheet, placed after "data" ------------------
$metaWs = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataWs)
$metaWs.Name = "metadata"

# --- Header row (bold, bordered, centered - matches the "data" sheet header) -
$metaWs.Range("B1").Value = "data_name"
$metaWs.Range("C1").Value = "data_id"
$metaWs.Range("D1").Value = "data_version"
$metaWs.Range("E1").Value = "data_version_created"
$metaWs.Range("F1").Value = "panel_query_time"
$metaWs.Range("G1").Value = "panel_get_request"

$headerRng = $metaWs.Range("B1:G1")
$headerRng.Font.Bold = $true
$headerRng.HorizontalAlignment = -4108
$headerRng.VerticalAlignment = -4160
$headerRng.Borders.LineStyle = 1
$headerRng.Borders.Weight = 2

# --- Data row -----------------------------------------------------------
$metaWs.Range("A2").Value = 0
$metaWs.Range("A2").HorizontalAlignment = -4108
$metaWs.Range("A2").VerticalAlignment = -4160
$metaWs.Range("A2").Borders.LineStyle = 1
$metaWs.Range("A2").Borders.Weight = 2

$metaWs.Range("B2").Value = "Hyperoxaluria"
$metaWs.Range("C2").Value = 119
$metaWs.Range("D2").Value = "'1.0"
$metaWs.Range("E2").Value = "2021-07-29T00:22:34.185526Z"
$metaWs.Range("F2").Value = "2021-10-05 14:34:05.680688"
$metaWs.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/119/?format=json"

# Keep "data" as the active sheet, matching the original workbook view state
$dataWs.Activate()
